$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row above the current row 3 ("Email profesional" row), which
# pushes rows 3-16 down to 4-17. Then drop the now-superfluous trailing
# row 17 so the sheet keeps its original 16-row extent.
# ---------------------------------------------------------------------------
$ws.Rows("3:3").Insert()
$ws.Rows("17:17").Delete()

# Re-apply the same cell formatting used by row 2 (Hosting/Hostinger) to the
# freshly inserted row 3, so it reuses the existing style ids instead of
# Excel fabricating brand-new ones.
$ws.Range("A2:I2").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122)
# Column F on this new row uses the USD currency format (same as the
# apify row), not the EUR format used elsewhere in column F.
$ws.Range("D5").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Populate the new row 3: "Hosting + backup diario" / Hostinger.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Hosting + backup diario"
$ws.Range("B3").Value = "Hostinger"
$ws.Range("C3").Value = "https://www.hostinger.com/"
$ws.Range("D3").Value = 928.37
$ws.Range("E3").Value = 0.21
$ws.Range("F3").Formula = "=D3*(1+E3)"
$ws.Range("G3").Value = "2 años"
$ws.Range("H3").Formula = "=F3/24"
$ws.Range("I3").Value = ""

# ---------------------------------------------------------------------------
# Rebuild the hyperlinks so they point at the shifted rows, and add the new
# one for the inserted row. Clearing via any single cell's Hyperlinks
# collection clears the whole sheet's hyperlink list in this host.
# ---------------------------------------------------------------------------
$ws.Range("C2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.hostinger.com/")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.hostinger.com/")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://apify.com/pricing/creator-plan")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://openrouter.ai/")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://groq.com/pricing")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.hostinger.com/")

# ---------------------------------------------------------------------------
# Cosmetic touch-ups that came along with the edit: column A grew wider to
# fit the new longer label, and the saved cursor position moved.
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 18.3
$ws.Range("C18").Select()
